$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.386.04'
$ws.Range('E2').Value = '  +3.20%  '

$ws.Range('D3').Value = '1.870.92'
$ws.Range('E3').Value = '  +1.50%  '

$ws.Range('E4').Value = '  -0.28%  '

$ws.Range('D5').Value = "'339.20"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.80%  '

$ws.Range('E6').Value = '  -0.30%  '

$ws.Range('D7').Value = "'0.4703"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.00%  '

$ws.Range('E8').Value = '  +3.40%  '

$ws.Range('D9').Value = "'47.41"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.17%  '

$ws.Range('D10').Value = "'0.08017"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.70%  '

$ws.Range('E11').Value = '  +2.36%  '

$ws.Range('D12').Value = "'21.86"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.07%  '

$ws.Range('D13').Value = '1.876.69'
$ws.Range('E13').Value = '  +1.70%  '

$ws.Range('D14').Value = "'5.988"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.52%  '

$ws.Range('D15').Value = "'7.241"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.97%  '

$ws.Range('D16').Value = "'91.22"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.82%  '

$ws.Range('E17').Value = '  -0.38%  '

$ws.Range('E18').Value = '  +1.17%  '

$ws.Range('D19').Value = "'0.06616"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.20%  '

$ws.Range('D20').Value = "'17.56"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.54%  '

$ws.Range('E21').Value = '  -0.34%  '

$ws.Range('D22').Value = '28.396.64'
$ws.Range('E22').Value = '  +3.25%  '

$ws.Range('D23').Value = "'5.452"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.94%  '

$ws.Range('D24').Value = "'11.04"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.44%  '

$ws.Range('E25').Value = '  -1.13%  '

$ws.Range('D26').Value = '2.098.75'
$ws.Range('E26').Value = '  +1.78%  '

$ws.Range('D27').Value = "'160.17"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.98%  '

$ws.Range('D28').Value = "'19.73"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.72%  '

$ws.Range('E29').Value = '  +2.65%  '

$ws.Range('E30').Value = '  +3.01%  '

$ws.Range('D31').Value = "'120.04"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.92%  '

$ws.Range('D32').Value = "'0.9691"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.15%  '

$ws.Range('D33').Value = "'0.09475"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.68%  '

$ws.Range('E34').Value = '  +0.03%  '

$ws.Range('D35').Value = "'1.376"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.06%  '

$ws.Range('E36').Value = '  +1.96%  '

$ws.Range('D37').Value = "'0.06082"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.25%  '

$ws.Range('D38').Value = "'0.02247"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.17%  '

$ws.Range('D39').Value = "'8.377"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.56%  '

$ws.Range('D40').Value = "'1.187"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.23%  '

$ws.Range('D41').Value = "'0.5946"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.90%  '

$ws.Range('E42').Value = '  -0.31%  '

$ws.Range('D43').Value = "'0.1870"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.42%  '

$ws.Range('E44').Value = '  +2.21%  '

$ws.Range('D45').Value = "'1.294"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.52%  '

$ws.Range('D46').Value = "'0.5583"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.27%  '

$ws.Range('D47').Value = "'12.18"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.95%  '

$ws.Range('D48').Value = "'1.954"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.52%  '

$ws.Range('D49').Value = "'0.06855"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.87%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').Value = "'2.043"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.15%  '

$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D51').Value = "'111.30"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.98%  '
